$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "tank"
$ws.Range("B29").Value = "agriculture"
$ws.Range("B36").Value = "settlement"
$ws.Range("B40").Value = "settlement"
$ws.Range("B42").Value = "crop"
$ws.Range("B47").Value = "river"
$ws.Range("B54").Value = "river"
$ws.Range("B58").Value = "tank"
$ws.Range("B80").Value = "river"
$ws.Range("B90").Value = "river"
$ws.Range("B94").Value = "road_n_railway"
$ws.Range("B98").Value = "crop"
$ws.Range("B103").Value = "river"
$ws.Range("B111").Value = "agriculture"
$ws.Range("B113").Value = "road_n_railway"
$ws.Range("B122").Value = "crop"
$ws.Range("B127").Value = "road_n_railway"
$ws.Range("B133").Value = "grassland"
$ws.Range("B137").Value = "forest"
$ws.Range("B141").Value = "river"
$ws.Range("B149").Value = "river"
$ws.Range("B150").Value = "river"
$ws.Range("B154").Value = "land_without_scrub"
$ws.Range("B160").Value = "river"
$ws.Range("B168").Value = "agriculture"
$ws.Range("B173").Value = "forest"
$ws.Range("B182").Value = "crop"
$ws.Range("B206").Value = "settlement"
$ws.Range("B207").Value = "settlement"
$ws.Range("B208").Value = "agriculture"
$ws.Range("B209").Value = "forest"
$ws.Range("B210").Value = "river"
$ws.Range("B219").Value = "forest"
$ws.Range("B223").Value = "land_without_scrub"
$ws.Range("B224").Value = "road_n_railway"
$ws.Range("B225").Value = "crop"
$ws.Range("B226").Value = "river"
$ws.Range("B228").Value = "river"
$ws.Range("B230").Value = "tank"
$ws.Range("B231").Value = "tank"
$ws.Range("B238").Value = "river"
$ws.Range("B247").Value = "crop"
$ws.Range("B250").Value = "road_n_railway"
$ws.Range("B251").Value = "crop"
$ws.Range("B268").Value = "tank"
$ws.Range("B273").Value = "road_n_railway"
$ws.Range("B274").Value = "river"
$ws.Range("B277").Value = "tank"
$ws.Range("B278").Value = "river"
$ws.Range("B279").Value = "river"
$ws.Range("B294").Value = "grassland"
$ws.Range("B304").Value = "settlement"
$ws.Range("B311").Value = "forest"
$ws.Range("B318").Value = "agriculture"
$ws.Range("B322").Value = "agriculture"
$ws.Range("B333").Value = "settlement"
$ws.Range("B335").Value = "settlement"
$ws.Range("B347").Value = "tank"
$ws.Range("B356").Value = "agriculture"
$ws.Range("B357").Value = "agriculture"
$ws.Range("B372").Value = "tank"
$ws.Range("B383").Value = "grassland"
$ws.Range("B393").Value = "crop"
$ws.Range("B398").Value = "land_without_scrub"
$ws.Range("B401").Value = "river"
$ws.Range("B410").Value = "agriculture"
$ws.Range("B414").Value = "agriculture"
$ws.Range("B417").Value = "grassland"
$ws.Range("B418").Value = "grassland"
$ws.Range("B421").Value = "agriculture"
$ws.Range("B423").Value = "grassland"
$ws.Range("B426").Value = "tank"
$ws.Range("B427").Value = "crop"
$ws.Range("B428").Value = "agriculture"
$ws.Range("B435").Value = "river"
$ws.Range("B439").Value = "river"
$ws.Range("B451").Value = "road_n_railway"
$ws.Range("B455").Value = "settlement"
$ws.Range("B459").Value = "agriculture"
$ws.Range("B460").Value = "forest"
$ws.Range("B462").Value = "crop"
$ws.Range("B467").Value = "crop"
$ws.Range("B474").Value = "agriculture"
$ws.Range("B475").Value = "road_n_railway"
$ws.Range("B504").Value = "tank"
